$wb = $excel.ActiveWorkbook

# Add one day of 2025 crime data (2025-11-29): update the "2025" (column L)
# totals on the Citywide Totals sheet, the By Neighborhood summary sheet, and
# every individual neighborhood sheet that had at least one incident that day.

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 6084
$ws.Range('L3').Value = 6624
$ws.Range('L4').Value = 1628
$ws.Range('L5').Value = 394
$ws.Range('L6').Value = 5426
$ws.Range('L7').Value = 20156

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 403
$ws.Range('L3').Value = 473
$ws.Range('L6').Value = 327
$ws.Range('L7').Value = 1338

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 138
$ws.Range('L3').Value = 179
$ws.Range('L7').Value = 440

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 245
$ws.Range('L3').Value = 319
$ws.Range('L4').Value = 61
$ws.Range('L6').Value = 260
$ws.Range('L7').Value = 907

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 233
$ws.Range('L3').Value = 273
$ws.Range('L7').Value = 771

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L6').Value = 157
$ws.Range('L7').Value = 651
$ws.Range('L8').Value = 1338
$ws.Range('L9').Value = 116
$ws.Range('L10').Value = 134
$ws.Range('L12').Value = 46
$ws.Range('L15').Value = 164
$ws.Range('L18').Value = 137
$ws.Range('L20').Value = 512
$ws.Range('L29').Value = 1129
$ws.Range('L31').Value = 198
$ws.Range('L33').Value = 907
$ws.Range('L34').Value = 114
$ws.Range('L36').Value = 255
$ws.Range('L37').Value = 771
$ws.Range('L42').Value = 645
$ws.Range('L49').Value = 109
$ws.Range('L51').Value = 253
$ws.Range('L52').Value = 429
$ws.Range('L54').Value = 438
$ws.Range('L55').Value = 213
$ws.Range('L60').Value = 130
$ws.Range('L63').Value = 57
$ws.Range('L67').Value = 699
$ws.Range('L78').Value = 264
$ws.Range('L83').Value = 440
$ws.Range('L85').Value = 998
$ws.Range('L86').Value = 130
$ws.Range('L90').Value = 210
$ws.Range('L96').Value = 224
$ws.Range('L97').Value = 163
$ws.Range('L101').Value = 20156

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L4').Value = 11
$ws.Range('L7').Value = 198

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 276
$ws.Range('L7').Value = 699

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('L4').Value = 15
$ws.Range('L7').Value = 109

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L4').Value = 36
$ws.Range('L7').Value = 438

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L3').Value = 436
$ws.Range('L6').Value = 276
$ws.Range('L7').Value = 1129

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L4').Value = 16
$ws.Range('L7').Value = 157

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 176
$ws.Range('L7').Value = 645

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L2').Value = 52
$ws.Range('L7').Value = 134

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L3').Value = 88
$ws.Range('L4').Value = 30
$ws.Range('L6').Value = 74
$ws.Range('L7').Value = 264

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 62
$ws.Range('L7').Value = 213

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 71
$ws.Range('L7').Value = 224

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L6').Value = 121
$ws.Range('L7').Value = 512

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L4').Value = 15
$ws.Range('L7').Value = 137

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L3').Value = 83
$ws.Range('L7').Value = 255

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 224
$ws.Range('L4').Value = 46
$ws.Range('L5').Value = 19
$ws.Range('L7').Value = 651

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L5').Value = 1
$ws.Range('L7').Value = 114

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L6').Value = 34
$ws.Range('L7').Value = 164

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('L2').Value = 36
$ws.Range('L7').Value = 116

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L2').Value = 40
$ws.Range('L7').Value = 163

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 69
$ws.Range('L7').Value = 130

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L6').Value = 62
$ws.Range('L7').Value = 210

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L4').Value = 37
$ws.Range('L7').Value = 253

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('L3').Value = 42
$ws.Range('L7').Value = 130

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 300
$ws.Range('L6').Value = 206
$ws.Range('L7').Value = 998

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 133
$ws.Range('L7').Value = 429

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('L2').Value = 14
$ws.Range('L7').Value = 46
